# Auto-generated edit script: updates crypto price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'34.484.26"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.Value = "'  +13.31%  "
$c.ClearFormats()
$c = $ws.Range("D3")
$c.Value = "'1.822.27"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.Value = "'  +7.95%  "
$c.ClearFormats()
$c = $ws.Range("E4")
$c.Value = "'  +0.49%  "
$c.ClearFormats()
$c = $ws.Range("D5")
$c.Value = "'232.67"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.Value = "'  +5.26%  "
$c.ClearFormats()
$c = $ws.Range("E6")
$c.Value = "'  +5.18%  "
$c.ClearFormats()
$c = $ws.Range("E7")
$c.Value = "'  +0.55%  "
$c.ClearFormats()
$c = $ws.Range("D8")
$c.Value = "'31.60"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.Value = "'  +4.11%  "
$c.ClearFormats()
$c = $ws.Range("D9")
$c.Value = "'46.36"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.Value = "'  +4.67%  "
$c.ClearFormats()
$c = $ws.Range("E10")
$c.Value = "'  +7.50%  "
$c.ClearFormats()
$c = $ws.Range("D11")
$c.Value = "'0.0684"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.Value = "'  +9.68%  "
$c.ClearFormats()
$c = $ws.Range("E12")
$c.Value = "'  +3.42%  "
$c.ClearFormats()
$c = $ws.Range("D13")
$c.Value = "'2.083.38"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.Value = "'  +8.05%  "
$c.ClearFormats()
$c = $ws.Range("D14")
$c.Value = "'1.818.47"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.Value = "'  +7.92%  "
$c.ClearFormats()
$c = $ws.Range("E15")
$c.Value = "'  +4.51%  "
$c.ClearFormats()
$c = $ws.Range("D16")
$c.Value = "'34.499.40"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.Value = "'  +13.48%  "
$c.ClearFormats()
$c = $ws.Range("E17")
$c.Value = "'  -4.58%  "
$c.ClearFormats()
$c = $ws.Range("E18")
$c.Value = "'  +8.48%  "
$c.ClearFormats()
$c = $ws.Range("D19")
$c.Value = "'71.05"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.Value = "'  +7.89%  "
$c.ClearFormats()
$c = $ws.Range("D20")
$c.Value = "'261.37"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.Value = "'  +6.02%  "
$c.ClearFormats()
$c = $ws.Range("D21")
$c.Value = "'0.0₃0755"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.Value = "'  +4.85%  "
$c.ClearFormats()
$c = $ws.Range("E22")
$c.Value = "'  +0.27%  "
$c.ClearFormats()
$c = $ws.Range("D23")
$c.Value = "'10.52"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.Value = "'  +3.00%  "
$c.ClearFormats()
$c = $ws.Range("E24")
$c.Value = "'  +2.72%  "
$c.ClearFormats()
$c = $ws.Range("D25")
$c.Value = "'2.20"
$c.ClearFormats()
$c = $ws.Range("D26")
$c.Value = "'162.11"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.Value = "'  +2.14%  "
$c.ClearFormats()
$c = $ws.Range("D27")
$c.Value = "'16.92"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.Value = "'  +6.41%  "
$c.ClearFormats()
$c = $ws.Range("E28")
$c.Value = "'  +5.08%  "
$c.ClearFormats()
$c = $ws.Range("D29")
$c.Value = "'7.16"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.Value = "'  +5.81%  "
$c.ClearFormats()
$c = $ws.Range("E30")
$c.Value = "'  +0.42%  "
$c.ClearFormats()
$c = $ws.Range("D31")
$c.Value = "'3.83"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.Value = "'  +9.31%  "
$c.ClearFormats()
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.Value = "'1.23"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.Value = "'  +7.46%  "
$c.ClearFormats()
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.Value = "'0.0516"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.Value = "'  +2.95%  "
$c.ClearFormats()
$c = $ws.Range("E34")
$c.Value = "'  +8.04%  "
$c.ClearFormats()
$c = $ws.Range("D35")
$c.Value = "'1.589.41"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.Value = "'  +5.17%  "
$c.ClearFormats()
$c = $ws.Range("E36")
$c.Value = "'  +5.93%  "
$c.ClearFormats()
$c = $ws.Range("D37")
$c.Value = "'1.06"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.Value = "'  +2.92%  "
$c.ClearFormats()
$c = $ws.Range("D38")
$c.Value = "'85.78"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.Value = "'  +8.26%  "
$c.ClearFormats()
$c = $ws.Range("D39")
$c.Value = "'0.0189"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.Value = "'  +5.14%  "
$c.ClearFormats()
$c = $ws.Range("E40")
$c.Value = "'  +7.91%  "
$c.ClearFormats()
$c = $ws.Range("E41")
$c.Value = "'  +1.86%  "
$c.ClearFormats()
$c = $ws.Range("E42")
$c.Value = "'  +2.59%  "
$c.ClearFormats()
$c = $ws.Range("D43")
$c.Value = "'0.918"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.Value = "'  +7.79%  "
$c.ClearFormats()
$c = $ws.Range("D44")
$c.Value = "'2.13"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.Value = "'  +6.39%  "
$c.ClearFormats()
$c = $ws.Range("D45")
$c.Value = "'0.0523"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.Value = "'  +3.75%  "
$c.ClearFormats()
$c = $ws.Range("E46")
$c.Value = "'  +6.69%  "
$c.ClearFormats()
$c = $ws.Range("D47")
$c.Value = "'1.975.83"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.Value = "'  +8.29%  "
$c.ClearFormats()
$c = $ws.Range("D48")
$c.Value = "'53.67"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.Value = "'  +2.99%  "
$c.ClearFormats()
$c = $ws.Range("D49")
$c.Value = "'5.74"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.Value = "'  +5.78%  "
$c.ClearFormats()
$c = $ws.Range("D50")
$c.Value = "'1.00"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.Value = "'  +0.34%  "
$c.ClearFormats()
$c = $ws.Range("D51")
$c.Value = "'0.0₆0123"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.Value = "'  +8.36%  "
$c.ClearFormats()
